# The deck ships two theme parts:
#   theme1.xml -> "Office Theme"  (currently only used by the Notes Master)
#   theme2.xml -> "Integral" / "Red Violet" colours (used by the one Slide Master,
#                 and therefore by every slide's active design)
#
# The authored edit swaps the two themes' contents: the design actually
# applied to the slides becomes the default "Office Theme" colour set,
# while the (until then unused-by-slides) "Integral" theme content moves
# to the other theme part.
#
# The only part of that swap an end user can drive through the PowerPoint
# object model is the *applied* theme's colour scheme (it's what
# Slide.ColorScheme / Slide.ThemeColorScheme exposes) - so we push the
# "Office Theme" RGB values onto the presentation's active colour scheme,
# which is what every slide (and the masters) actually render with.

$p = $ppt.ActivePresentation

function Convert-HexToRGB($hex) {
    # PowerPoint RGB longs are packed as 0x00BBGGRR (same order COM returns).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# "Office Theme" colour scheme, in the standard 12-slot theme colour order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide = $p.Slides.Item(1)
$colorScheme = $slide.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Convert-HexToRGB $officeThemeColors[$i - 1]
}
